$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of match data to append after the existing last row (36).
$rows = @(
    @{
        idx = 36
        country = "kuwait"
        tournament = "premier-league"
        season = "2023-2024"
        matchDate = 45255.63888888889
        home = "Al Jahra"
        homeGoals = 0
        away = "Al Naser"
        awayGoals = 3
        homeOpenOdds = 3.19
        homeOpenDt = "25/11/2023 03:42"
        homeCloseOdds = 3.56
        homeCloseDt = "25/11/2023 13:07"
        drawOpenOdds = 3.43
        drawOpenDt = "25/11/2023 03:42"
        drawCloseOdds = 3.47
        drawCloseDt = "25/11/2023 13:25"
        awayOpenOdds = 2.06
        awayOpenDt = "25/11/2023 03:42"
        awayCloseOdds = 1.92
        awayCloseDt = "25/11/2023 13:16"
        url = "https://www.betexplorer.com/football/kuwait/premier-league/al-jahra-al-naser/G6wuXpAb/"
    },
    @{
        idx = 37
        country = "kuwait"
        tournament = "premier-league"
        season = "2023-2024"
        matchDate = 45255.75694444445
        home = "Al-Fahaheel"
        homeGoals = 3
        away = "Kazma SC"
        awayGoals = 4
        homeOpenOdds = 3.42
        homeOpenDt = "25/11/2023 06:12"
        homeCloseOdds = 3.62
        homeCloseDt = "25/11/2023 18:09"
        drawOpenOdds = 3.65
        drawOpenDt = "25/11/2023 06:12"
        drawCloseOdds = 3.84
        drawCloseDt = "25/11/2023 18:09"
        awayOpenOdds = 1.9
        awayOpenDt = "25/11/2023 06:12"
        awayCloseOdds = 1.82
        awayCloseDt = "25/11/2023 18:09"
        url = "https://www.betexplorer.com/football/kuwait/premier-league/al-fahaheel-kazma-sc/xGvyYQfh/"
    },
    @{
        idx = 38
        country = "kuwait"
        tournament = "premier-league"
        season = "2023-2024"
        matchDate = 45256.63888888889
        home = "Khaitan"
        homeGoals = 1
        away = "Al Shabab"
        awayGoals = 4
        homeOpenOdds = 2.81
        homeOpenDt = "26/11/2023 03:43"
        homeCloseOdds = 3.42
        homeCloseDt = "26/11/2023 15:03"
        drawOpenOdds = 3.39
        drawOpenDt = "26/11/2023 03:43"
        drawCloseOdds = 3.41
        drawCloseDt = "26/11/2023 15:03"
        awayOpenOdds = 2.22
        awayOpenDt = "26/11/2023 03:43"
        awayCloseOdds = 2
        awayCloseDt = "26/11/2023 15:03"
        url = "https://www.betexplorer.com/football/kuwait/premier-league/khaitan-al-shabab/UkyTZnQu/"
    }
)

$lastExistingRow = 36
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy formatting (font/border/alignment/number formats) from the last
    # pre-existing data row so the new row matches the sheet's look (bold
    # bordered index column, date-formatted match-date column, etc.).
    $ws.Range("A" + $lastExistingRow + ":V" + $lastExistingRow).Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data.idx
    $ws.Cells.Item($r, 2).Value = $data.country
    $ws.Cells.Item($r, 3).Value = $data.tournament
    $ws.Cells.Item($r, 4).Value = $data.season
    $ws.Cells.Item($r, 5).Value = $data.matchDate
    $ws.Cells.Item($r, 6).Value = $data.home
    $ws.Cells.Item($r, 7).Value = $data.homeGoals
    $ws.Cells.Item($r, 8).Value = $data.away
    $ws.Cells.Item($r, 9).Value = $data.awayGoals
    $ws.Cells.Item($r, 10).Value = $data.homeOpenOdds
    $ws.Cells.Item($r, 11).Value = $data.homeOpenDt
    $ws.Cells.Item($r, 12).Value = $data.homeCloseOdds
    $ws.Cells.Item($r, 13).Value = $data.homeCloseDt
    $ws.Cells.Item($r, 14).Value = $data.drawOpenOdds
    $ws.Cells.Item($r, 15).Value = $data.drawOpenDt
    $ws.Cells.Item($r, 16).Value = $data.drawCloseOdds
    $ws.Cells.Item($r, 17).Value = $data.drawCloseDt
    $ws.Cells.Item($r, 18).Value = $data.awayOpenOdds
    $ws.Cells.Item($r, 19).Value = $data.awayOpenDt
    $ws.Cells.Item($r, 20).Value = $data.awayCloseOdds
    $ws.Cells.Item($r, 21).Value = $data.awayCloseDt
    $ws.Cells.Item($r, 22).Value = $data.url
}
